$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Locate the paragraph that currently reads:
#      <ab><margin>left-top</margin>
#    (w14:paraId="0000000E") and replace it with two paragraphs:
#      <ab><margin>left-middle</margin>
#      <render>tall</render>
#    while also stripping the pBdr/shd noise from pPr and giving the
#    paragraph mark / "left-middle" run the updated run properties
#    described by the diff.
# ------------------------------------------------------------------

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*<ab><margin>left-top</margin>*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the '<ab><margin>left-top</margin>' paragraph"
}

$newParasXml = '<w:p w:rsidR="00000000" w:rsidDel="00000000" w:rsidP="00000000" w:rsidRDefault="00000000" w:rsidRPr="00000000" w14:paraId="0000000E"><w:pPr><w:widowControl w:val="0"/><w:contextualSpacing w:val="0"/><w:rPr><w:rFonts w:ascii="Courier New" w:cs="Courier New" w:eastAsia="Courier New" w:hAnsi="Courier New"/><w:color w:val="7f6000"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r w:rsidDel="00000000" w:rsidR="00000000" w:rsidRPr="00000000"><w:rPr><w:rFonts w:ascii="Courier New" w:cs="Courier New" w:eastAsia="Courier New" w:hAnsi="Courier New"/><w:color w:val="7f6000"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">&lt;ab&gt;&lt;margin&gt;</w:t></w:r><w:r w:rsidDel="00000000" w:rsidR="00000000" w:rsidRPr="00000000"><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">left-middle</w:t></w:r><w:r w:rsidDel="00000000" w:rsidR="00000000" w:rsidRPr="00000000"><w:rPr><w:rFonts w:ascii="Courier New" w:cs="Courier New" w:eastAsia="Courier New" w:hAnsi="Courier New"/><w:color w:val="7f6000"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">&lt;/margin&gt;</w:t></w:r></w:p><w:p w:rsidR="00000000" w:rsidDel="00000000" w:rsidP="00000000" w:rsidRDefault="00000000" w:rsidRPr="00000000" w14:paraId="0000000E"><w:pPr><w:widowControl w:val="0"/><w:contextualSpacing w:val="0"/><w:rPr><w:rFonts w:ascii="Courier New" w:cs="Courier New" w:eastAsia="Courier New" w:hAnsi="Courier New"/><w:color w:val="7f6000"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r w:rsidDel="00000000" w:rsidR="00000000" w:rsidRPr="00000000"><w:rPr><w:rFonts w:ascii="Courier New" w:cs="Courier New" w:eastAsia="Courier New" w:hAnsi="Courier New"/><w:color w:val="7f6000"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">&lt;render&gt;tall&lt;/render&gt;</w:t></w:r></w:p>'

$pkgXml = '<?xml version="1.0"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>' + $newParasXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

[void]$target.Range.InsertXML($pkgXml)

# ------------------------------------------------------------------
# 2) Add a footer distance (w:pgMar/@w:footer="720", i.e. 36pt) to
#    the section properties.
# ------------------------------------------------------------------
$d.PageSetup.FooterDistance = 36
